$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Palo Alto, CA" / "Stanford" case-study header block (rows 51-53).
# The hotel listing rows below it (55-59) stay put, so we clear the cell
# contents in place rather than deleting whole rows (which would shift
# everything beneath upward).
$ws.Range("A51:D53").ClearContents()

# Add the new "Westfield NJ" case-study details that were appended below
# the existing data.
$ws.Range("A76").Value = "westfield Inn, BW Signature"
$ws.Range("B76").Value = 41
$ws.Range("A78").Value = "No major college other than Union Co. Community College"
$ws.Range("A79").Value = "No inpatient hospital "
$ws.Range("A82").Value = "Clarendon VA"

# Match the saved selection / scroll position.
$ws.Range("A82").Select()
$excel.ActiveWindow.ScrollRow = 75
$excel.ActiveWindow.ScrollColumn = 1
